# Add article 51 ("Healing the fabrics") as a new row at the bottom of Sheet1,
# mirroring the layout of the existing rows (Ser, Date, Ayats, Content, Author, Tags).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52

$tags = @'
Family to friends to enemies, The skin of a human, The greatest warrior, Love of wife, Forgiveness
'@

$body = @'
h1: Healing the fabrics
p.note: This topic can make our lives better. Please brothers and sisters listen..
img.width-m-half: forts.jpg
p.caption.text-center: Photo by <a href=https://unsplash.com/@alirezapasandeh?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText target=_blank>Alireza Pasandeh</a> on Unsplash
p: There are small forts in our lives, like the one above. Each fort has walls around it, the anchors and the oil trench with spear headings protecting it from offensives. Some forts are built higher and some are built weaker. The forts that are bigger in size, control a little more of the land around them. They govern some more rules than the smaller forts.<b> Bigger forts do not break easily. They are taken down only in wars. They lead the trust of the people. </b>
p: Our families, friends and relationships are these forts. Each husband and wife is an institution residing in this fort. The brother, might be the most miserable person on the earth, is your line of defence against your bigger enemies. The sisters, who might do not speak their entire lives, are the best fighters to protect our institutions. The friends, who may criticise you at times, are your friends to protect your institution. Each person we interact, creates a pulse of emotions, memories and sparks that can protect only us. 
p: I am talking selfish, totally focused on my own success. The success of my self, is only possible after I build a strong fort with good roads, better rules and long lasting relations. I can not succeed without the existence of an institution. I can not walk for my entire life without the help of my wife, my brother and my cousin. Each relationship makes me a success. Every single bond makes me human. I can not talk loud to anybody, if I am so in love with my own success. The moment, I break a laid down principle, I react in anger, I make this crack in the fort that takes significant effort to repair. 
h3: But some people are really bad
p: Yes some hurt us more than others. Some even beat us to the ground. Some bully us in public. Some throw dust at us. Some are evil. Some are dancing in the mosques. Yes there are people who make us angry and are really bad. 
p.note: How do I say it…
p: They are the people who have a human underneath their skin. Some serial killers who are stuck in prisons, when they cry, they bring you down with them. The tears falling from great chests and strong shoulders are the heaviest. Some who killed our brothers and sisters on spearheads. They murdered them in public. These killers became the best fighters and commanders in wars. 
p: These same murderers and killers when raised their voices, they shook into the hearts of the enemies. They shook this earth. Their voice moved hearts. Their wives loved them. Their people followed them into destruction. They died leading a battalion of 400 people and all 400 of their people died with them. <b>Ikram (ra) after becoming a war criminal, became a commander in Battle of Yarmouk</b>. He died like no body of us can die. His voice was the loudest in the crowd. The moment he stood his ground, enemy could not stand with all its might. All 70,000 were down. We won Jerusalem and we raised higher.
img.width-m-half: yarmouk.jpg
p.caption.text-center: Ikrama (ra) on 4th Day of battle in Yarmouk, snippet from <a href=https://www.amazon.com/Sword-Allah-I-Akram/dp/0195977149 target=_blank>The Sword of Allah</a>
p: These murderers have a human underneath their skin. The worst attackers have the human part. Look at the emotions they carry. The tears they drop. The love they give. Look at how they protect you when you are falling. 
h3: Conclusion
p: Each relationship is a meaningful relationship. The person meeting us on the road to the wife inside our home. All these people are independent forts. We need to build stronger forts and harder institutions that we rise higher in this world. We can not do it alone, we need help of these other forts around us. We need to heal the roads leading to these other cities and towns. These fabrics woven into our hearts, needs to warmth. We need to forgive others. The forgiveness that puts a dent into their hearts. The dent that lasts an eternity and leads all of us to the success.
p: It is not possible, we heal these fabrics right away. It takes time. Give your relationships everything you can that they protect you. Only you will be benefitted. Be selfish. So that, only <b>you</b> succeed.
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.
'@

$ayats = @'
Surah Al Nisa, 101 - 106
'@

$ws.Cells.Item($newRow, 1).Value = 51
$ws.Cells.Item($newRow, 2).Value = 43881
$ws.Cells.Item($newRow, 2).NumberFormat = "d-mmm-yy"
$ws.Cells.Item($newRow, 3).Value = $ayats
$ws.Cells.Item($newRow, 4).Value = $body
$ws.Cells.Item($newRow, 5).Value = "Qasim Ali"
$ws.Cells.Item($newRow, 6).Value = $tags

$null = $ws.Range("C52").Select()
